$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

function Copy-CellFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# ---------------------------------------------------------------------------
# 1. Title cells (A1:B1) pick up the "no side borders" look already used by
#    I15 (top/bottom border only, regular font, left aligned).
# ---------------------------------------------------------------------------
Copy-CellFormat "I15" "A1"
Copy-CellFormat "I15" "B1"

# ---------------------------------------------------------------------------
# 2. New observation rows (19-21), added below the existing data.
# ---------------------------------------------------------------------------

# --- Row 19 : RASY @ Sainte-Louise (Chaudiere-Appalaches) ---
$ws.Range("A19").Value = 45041
$ws.Range("B19").Value = "RASY"
$ws.Range("C19").Value = "N/A"
$ws.Range("D19").Value = "Sainte-Louise"
$ws.Range("E19").Value = "Chaudière-Appalaches"
$ws.Range("F19").Value = "B"
$ws.Range("G19").Value = "Cote 1"
$ws.Range("H19").Value = "iNaturalist (https://www.inaturalist.org/observations/156684795)"
$ws.Range("I19").Value = "Julien Savoie"

Copy-CellFormat "A18" "A19"
Copy-CellFormat "B18" "B19"
Copy-CellFormat "C18" "C19"
Copy-CellFormat "D18" "D19"
Copy-CellFormat "E18" "E19"
Copy-CellFormat "F16" "F19"
Copy-CellFormat "G18" "G19"
Copy-CellFormat "H18" "H19"
Copy-CellFormat "I18" "I19"

# --- Row 20 : BUAM @ Rosemere (Laurentides) ---
$ws.Range("A20").Value = 45041
$ws.Range("B20").Value = "BUAM"
$ws.Range("C20").Value = "N/A"
$ws.Range("D20").Value = "Rosemère"
$ws.Range("E20").Value = "Laurentides"
$ws.Range("F20").Value = "A"
$ws.Range("G20").Value = "Cote 1"
$ws.Range("H20").Value = "iNaturalist (https://www.inaturalist.org/observations/156678026)"
$ws.Range("I20").Value = "ramenramen"

Copy-CellFormat "A18" "A20"
Copy-CellFormat "B10" "B20"
Copy-CellFormat "C18" "C20"
Copy-CellFormat "D18" "D20"
Copy-CellFormat "E18" "E20"
Copy-CellFormat "F18" "F20"
Copy-CellFormat "G18" "G20"
Copy-CellFormat "H18" "H20"
Copy-CellFormat "I18" "I20"

# --- Row 21 : BUAM @ Rosemere (Laurentides), second contact ---
$ws.Range("A21").Value = 45041
$ws.Range("B21").Value = "BUAM"
$ws.Range("C21").Value = "N/A"
$ws.Range("D21").Value = "Rosemère"
$ws.Range("E21").Value = "Laurentides"
$ws.Range("F21").Value = "A"
$ws.Range("G21").Value = "Cote 1"
$ws.Range("H21").Value = "iNaturalist (https://www.inaturalist.org/observations/156668393)"
$ws.Range("I21").Value = "philomene222"

Copy-CellFormat "A18" "A21"
Copy-CellFormat "B10" "B21"
Copy-CellFormat "C18" "C21"
Copy-CellFormat "D18" "D21"
Copy-CellFormat "E18" "E21"
Copy-CellFormat "F18" "F21"
Copy-CellFormat "G18" "G21"
Copy-CellFormat "H18" "H21"
Copy-CellFormat "I18" "I21"

# ---------------------------------------------------------------------------
# 3. Restore the selection the workbook was left with.
# ---------------------------------------------------------------------------
$ws.Range("F26").Select() | Out-Null
